# edit.ps1
# Applies the "Update July 2020 ppts" changes:
#   1. Remove the last two slides (sldId 260 / 261 -> slides/slide4.xml,
#      slides/slide5.xml) from the deck.
#   2. Refresh the cached text of the auto-updating "Date Placeholder"
#      field (type="datetimeFigureOut") on the slide master and on every
#      slide layout from 13/06/2020 to 05/07/2020.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Drop the trailing two slides (originally slide4.xml / slide5.xml,
#    the two full-bleed JPEG picture slides).
# ---------------------------------------------------------------------
for ($i = $p.Slides.Count; $i -ge 4; $i--) {
    $p.Slides.Item($i).Delete()
}

# ---------------------------------------------------------------------
# 2) Update the "Date Placeholder" shape text wherever it appears:
#    the slide master plus all slide layouts.
# ---------------------------------------------------------------------
function Update-DatePlaceholder($shapes) {
    $count = $shapes.Count
    for ($i = 1; $i -le $count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            if ($shp.TextFrame.TextRange.Text -eq "13/06/2020") {
                $shp.TextFrame.TextRange.Text = "05/07/2020"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}
